$p = $ppt.ActivePresentation

# --- 1) Update the cached "datetimeFigureOut" footer date from 3/22/2022 to
#        3/29/2022 everywhere it appears: the slide master and every one of
#        its custom (slide) layouts expose a "Date Placeholder" shape whose
#        text caches the date field's rendered value. ---
$newDate = "3/29/2022"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $layout = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Reposition the picture on slide 12 (the sweep/plot screenshot) up
#        and to the left, keeping its size the same. ---
$s12 = $p.Slides.Item(12)
for ($i = 1; $i -le $s12.Shapes.Count; $i++) {
    $sh = $s12.Shapes.Item($i)
    if ($sh.Name -like "Content Placeholder*") {
        $sh.Left = 111.72338582677165
        $sh.Top = 44.93748031496063
    }
}
